$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-05 Thursday" "2026-02-06 Friday"

Replace-Text "50÷2=25, 0" "54÷2=27, 0"
Replace-Text "30÷5=6, 0" "87÷8=10, 7"
Replace-Text "14÷5=2, 4" "69÷6=11, 3"
Replace-Text "94÷8=11, 6" "84÷9=9, 3"
Replace-Text "55÷8=6, 7" "48÷4=12, 0"

Replace-Text "15÷8=1, 7" "13÷5=2, 3"
Replace-Text "58÷5=11, 3" "96÷3=32, 0"
Replace-Text "59÷4=14, 3" "50÷8=6, 2"
Replace-Text "68÷8=8, 4" "87÷9=9, 6"
Replace-Text "81÷7=11, 4" "88÷6=14, 4"

Replace-Text "69÷2=34, 1" "20÷3=6, 2"
Replace-Text "69÷7=9, 6" "51÷2=25, 1"
Replace-Text "96÷2=48, 0" "15÷3=5, 0"
Replace-Text "23÷9=2, 5" "42÷3=14, 0"
Replace-Text "29÷8=3, 5" "45÷8=5, 5"

Replace-Text "12÷8=1, 4" "66÷4=16, 2"
Replace-Text "23÷8=2, 7" "85÷9=9, 4"
Replace-Text "34÷3=11, 1" "79÷2=39, 1"
Replace-Text "11÷4=2, 3" "46÷3=15, 1"
Replace-Text "31÷5=6, 1" "80÷3=26, 2"

Replace-Text "99÷8=12, 3" "26÷2=13, 0"
Replace-Text "19÷5=3, 4" "12÷5=2, 2"
Replace-Text "92÷5=18, 2" "63÷3=21, 0"
Replace-Text "50÷3=16, 2" "88÷7=12, 4"
Replace-Text "65÷3=21, 2" "50÷2=25, 0"

Write-Output "Done applying replacements."
